# Split documentation to have API specific section. Clarified Allergies.
# - bump Version / Date metadata values
# - mark the existing wrap/vertical-top alignment as explicitly "applied"
#   (Excel sets applyAlignment="true" on the cell styles once the
#   alignment is (re)applied through the UI/object model)

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "0.1.16-beta"
$meta.Range("B8").Value = "2023-06-13T11:38:47-05:00"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count

    $header = $ws.Range("A1:B1")
    $header.VerticalAlignment = -4160
    $header.WrapText = $true

    if ($rows -gt 1) {
        $body = $ws.Range("A2:B" + $rows)
        $body.VerticalAlignment = -4160
        $body.WrapText = $true
    }
}
